$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L6").Value = 0.8333333333333334
$ws.Range("M6").Value = 0.7692307692307693
$ws.Range("N6").Value = 0.8
$ws.Range("O6").Value = 0.010752688172043
$ws.Range("P6").Value = 0.1153846153846154
$ws.Range("Q6").Value = 0.0714285714285714
$ws.Range("R6").Value = 0.01307189542483659
$ws.Range("S6").Value = 0.1764705882352942
$ws.Range("T6").Value = 0.09803921568627445
$ws.Range("L9").Value = 0.8253968253968254
$ws.Range("M9").Value = 0.6666666666666666
$ws.Range("N9").Value = 0.7375886524822695
$ws.Range("O9").Value = 0.002816180235534982
$ws.Range("P9").Value = 0.01282051282051277
$ws.Range("Q9").Value = 0.00901722391084081
$ws.Range("R9").Value = 0.003423591658885664
$ws.Range("S9").Value = 0.01960784313725483
$ws.Range("T9").Value = 0.01237658183840895
$ws.Range("L11").Value = 0.8048780487804879
$ws.Range("M11").Value = 0.8461538461538461
$ws.Range("N11").Value = 0.8250000000000001
$ws.Range("O11").Value = -0.01770259638080252
$ws.Range("P11").Value = 0.1923076923076923
$ws.Range("Q11").Value = 0.09642857142857142
$ws.Range("R11").Value = -0.02152080344332855
$ws.Range("S11").Value = 0.2941176470588235
$ws.Range("T11").Value = 0.1323529411764706
$ws.Range("L16").Value = 0.8333333333333334
$ws.Range("M16").Value = 0.7692307692307693
$ws.Range("N16").Value = 0.8
$ws.Range("O16").Value = 0.010752688172043
$ws.Range("P16").Value = 0.1153846153846154
$ws.Range("Q16").Value = 0.0714285714285714
$ws.Range("R16").Value = 0.01307189542483659
$ws.Range("S16").Value = 0.1764705882352942
$ws.Range("T16").Value = 0.09803921568627445
$ws.Range("L19").Value = 0.8253968253968254
$ws.Range("M19").Value = 0.6666666666666666
$ws.Range("N19").Value = 0.7375886524822695
$ws.Range("O19").Value = 0.002816180235534982
$ws.Range("P19").Value = 0.01282051282051277
$ws.Range("Q19").Value = 0.00901722391084081
$ws.Range("R19").Value = 0.003423591658885664
$ws.Range("S19").Value = 0.01960784313725483
$ws.Range("T19").Value = 0.01237658183840895
$ws.Range("L21").Value = 0.8048780487804879
$ws.Range("M21").Value = 0.8461538461538461
$ws.Range("N21").Value = 0.8250000000000001
$ws.Range("O21").Value = -0.01770259638080252
$ws.Range("P21").Value = 0.1923076923076923
$ws.Range("Q21").Value = 0.09642857142857142
$ws.Range("R21").Value = -0.02152080344332855
$ws.Range("S21").Value = 0.2941176470588235
$ws.Range("T21").Value = 0.1323529411764706
$ws.Range("L25").Value = 0.9375
$ws.Range("M25").Value = 0.1923076923076923
$ws.Range("N25").Value = 0.3191489361702128
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = [double]"2.775557561562891E-17"
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = 0
$ws.Range("S25").Value = [double]"1.443289932012704E-16"
$ws.Range("T25").Value = 0
$ws.Range("L26").Value = 0.9375
$ws.Range("M26").Value = 0.1923076923076923
$ws.Range("N26").Value = 0.3191489361702128
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = [double]"2.775557561562891E-17"
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = 0
$ws.Range("S26").Value = [double]"1.443289932012704E-16"
$ws.Range("T26").Value = 0
$ws.Range("L29").Value = 0.9411764705882353
$ws.Range("M29").Value = 0.2051282051282051
$ws.Range("N29").Value = 0.3368421052631579
$ws.Range("O29").Value = 0.003676470588235281
$ws.Range("P29").Value = 0.01282051282051283
$ws.Range("Q29").Value = 0.0176931690929451
$ws.Range("R29").Value = 0.003921568627450966
$ws.Range("S29").Value = 0.06666666666666672
$ws.Range("T29").Value = 0.05543859649122799
$ws.Range("L36").Value = 0.6530612244897959
$ws.Range("M36").Value = 0.8205128205128205
$ws.Range("N36").Value = 0.7272727272727272
$ws.Range("O36").Value = 0.4030612244897959
$ws.Range("P36").Value = 0.8076923076923077
$ws.Range("Q36").Value = 0.7028824833702881
$ws.Range("R36").Value = 1.612244897959183
$ws.Range("S36").Value = 63.00000000000011
$ws.Range("T36").Value = 28.81818181818184
